$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 13.739149
$ws.Range("H2").Value = 41.217447
$ws.Range("I2").Value = 0.6130043224686931
$ws.Range("J2").Value = 0.6130043224686931
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 116.9511416666667
$ws.Range("N2").Value = 350.853425
$ws.Range("O2").Value = 0.411881549221027
$ws.Range("P2").Value = 0.411881549221027
$ws.Range("Q2").Value = 1606.809161078442
$ws.Range("R2").Value = 14461.28244970598
$ws.Range("S2").Value = 0.2524851700175913
$ws.Range("T2").Value = 0.2524851700175913

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 13.739149
$ws.Range("H3").Value = 41.217447
$ws.Range("I3").Value = 0.6130043224686931
$ws.Range("J3").Value = 0.6130043224686931
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 133.0753813333333
$ws.Range("N3").Value = 399.226144
$ws.Range("O3").Value = 0.468668312644395
$ws.Range("P3").Value = 0.468668312644395
$ws.Range("Q3").Value = 1828.342492370485
$ws.Range("R3").Value = 16455.08243133437
$ws.Range("S3").Value = 0.287295701455123
$ws.Range("T3").Value = 0.287295701455123

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 13.739149
$ws.Range("H4").Value = 41.217447
$ws.Range("I4").Value = 0.6130043224686931
$ws.Range("J4").Value = 0.6130043224686931
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 33.91710566666666
$ws.Range("N4").Value = 101.751317
$ws.Range("O4").Value = 0.119450138134578
$ws.Range("P4").Value = 0.119450138134578
$ws.Range("Q4").Value = 465.9921684030776
$ws.Range("R4").Value = 4193.929515627699
$ws.Range("S4").Value = 0.07322345099597881
$ws.Range("T4").Value = 0.07322345099597881

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.034036666666666
$ws.Range("H5").Value = 12.10211
$ws.Range("I5").Value = 0.1799879973398545
$ws.Range("J5").Value = 0.1799879973398545
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 116.9511416666667
$ws.Range("N5").Value = 350.853425
$ws.Range("O5").Value = 0.411881549221027
$ws.Range("P5").Value = 0.411881549221027
$ws.Range("Q5").Value = 471.7851936918611
$ws.Range("R5").Value = 4246.06674322675
$ws.Range("S5").Value = 0.07413373518552936
$ws.Range("T5").Value = 0.07413373518552936

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.034036666666666
$ws.Range("H6").Value = 12.10211
$ws.Range("I6").Value = 0.1799879973398545
$ws.Range("J6").Value = 0.1799879973398545
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 133.0753813333333
$ws.Range("N6").Value = 399.226144
$ws.Range("O6").Value = 0.468668312644395
$ws.Range("P6").Value = 0.468668312644395
$ws.Range("Q6").Value = 536.8309677293155
$ws.Range("R6").Value = 4831.478709563839
$ws.Range("S6").Value = 0.08435467100951348
$ws.Range("T6").Value = 0.08435467100951348

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.034036666666666
$ws.Range("H7").Value = 12.10211
$ws.Range("I7").Value = 0.1799879973398545
$ws.Range("J7").Value = 0.1799879973398545
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 33.91710566666666
$ws.Range("N7").Value = 101.751317
$ws.Range("O7").Value = 0.119450138134578
$ws.Range("P7").Value = 0.119450138134578
$ws.Range("Q7").Value = 136.8228478865411
$ws.Range("R7").Value = 1231.40563097887
$ws.Range("S7").Value = 0.02149959114481169
$ws.Range("T7").Value = 0.02149959114481169

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.639623666666666
$ws.Range("H8").Value = 13.918871
$ws.Range("I8").Value = 0.2070076801914524
$ws.Range("J8").Value = 0.2070076801914524
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 116.9511416666667
$ws.Range("N8").Value = 350.853425
$ws.Range("O8").Value = 0.411881549221027
$ws.Range("P8").Value = 0.411881549221027
$ws.Range("Q8").Value = 542.6092847203528
$ws.Range("R8").Value = 4883.483562483175
$ws.Range("S8").Value = 0.08526264401790631
$ws.Range("T8").Value = 0.08526264401790631

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.639623666666666
$ws.Range("H9").Value = 13.918871
$ws.Range("I9").Value = 0.2070076801914524
$ws.Range("J9").Value = 0.2070076801914524
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 133.0753813333333
$ws.Range("N9").Value = 399.226144
$ws.Range("O9").Value = 0.468668312644395
$ws.Range("P9").Value = 0.468668312644395
$ws.Range("Q9").Value = 617.4196886848248
$ws.Range("R9").Value = 5556.777198163423
$ws.Range("S9").Value = 0.09701794017975855
$ws.Range("T9").Value = 0.09701794017975855

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.639623666666666
$ws.Range("H10").Value = 13.918871
$ws.Range("I10").Value = 0.2070076801914524
$ws.Range("J10").Value = 0.2070076801914524
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 33.91710566666666
$ws.Range("N10").Value = 101.751317
$ws.Range("O10").Value = 0.119450138134578
$ws.Range("P10").Value = 0.119450138134578
$ws.Range("Q10").Value = 157.3626061559007
$ws.Range("R10").Value = 1416.263455403107
$ws.Range("S10").Value = 0.02472709599378755
$ws.Range("T10").Value = 0.02472709599378755
